$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7894570000000001
$ws.Range("N2").Value = 2.368371
$ws.Range("O2").Value = 0.09838606084581891
$ws.Range("P2").Value = 0.09838606084581894
$ws.Range("Q2").Value = 0.5898183243830001
$ws.Range("R2").Value = 5.308364919447
$ws.Range("S2").Value = 0.003046971887586011
$ws.Range("T2").Value = 0.003046971887586012
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.625751333333334
$ws.Range("N3").Value = 16.877254
$ws.Range("O3").Value = 0.7011091332204036
$ws.Range("P3").Value = 0.7011091332204038
$ws.Range("Q3").Value = 4.203105710408667
$ws.Range("R3").Value = 37.827951393678
$ws.Range("S3").Value = 0.02171303333711169
$ws.Range("T3").Value = 0.02171303333711169
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.608865333333333
$ws.Range("N4").Value = 4.826596
$ws.Range("O4").Value = 0.2005048059337773
$ws.Range("P4").Value = 0.2005048059337774
$ws.Range("Q4").Value = 1.202013858974667
$ws.Range("R4").Value = 10.818124730772
$ws.Range("S4").Value = 0.006209543321014777
$ws.Range("T4").Value = 0.006209543321014779
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7894570000000001
$ws.Range("N5").Value = 2.368371
$ws.Range("O5").Value = 0.09838606084581891
$ws.Range("P5").Value = 0.09838606084581894
$ws.Range("Q5").Value = 15.58876976089634
$ws.Range("R5").Value = 140.298927848067
$ws.Range("S5").Value = 0.08053080289289
$ws.Range("T5").Value = 0.08053080289289002
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.625751333333334
$ws.Range("N6").Value = 16.877254
$ws.Range("O6").Value = 0.7011091332204036
$ws.Range("P6").Value = 0.7011091332204038
$ws.Range("Q6").Value = 111.0871678475065
$ws.Range("R6").Value = 999.7845106275581
$ws.Range("S6").Value = 0.5738707386837785
$ws.Range("T6").Value = 0.5738707386837787
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.608865333333333
$ws.Range("N7").Value = 4.826596
$ws.Range("O7").Value = 0.2005048059337773
$ws.Range("P7").Value = 0.2005048059337774
$ws.Range("Q7").Value = 31.76896431043245
$ws.Range("R7").Value = 285.920678793892
$ws.Range("S7").Value = 0.1641168765871611
$ws.Range("T7").Value = 0.1641168765871611
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7894570000000001
$ws.Range("N8").Value = 2.368371
$ws.Range("O8").Value = 0.09838606084581891
$ws.Range("P8").Value = 0.09838606084581894
$ws.Range("Q8").Value = 2.866517577543
$ws.Range("R8").Value = 25.798658197887
$ws.Range("S8").Value = 0.01480828606534289
$ws.Range("T8").Value = 0.0148082860653429
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.625751333333334
$ws.Range("N9").Value = 16.877254
$ws.Range("O9").Value = 0.7011091332204036
$ws.Range("P9").Value = 0.7011091332204038
$ws.Range("Q9").Value = 20.427097465582
$ws.Range("R9").Value = 183.843877190238
$ws.Range("S9").Value = 0.1055253611995133
$ws.Range("T9").Value = 0.1055253611995134
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.608865333333333
$ws.Range("N10").Value = 4.826596
$ws.Range("O10").Value = 0.2005048059337773
$ws.Range("P10").Value = 0.2005048059337774
$ws.Range("Q10").Value = 5.841788416468
$ws.Range("R10").Value = 52.576095748212
$ws.Range("S10").Value = 0.03017838602560145
$ws.Range("T10").Value = 0.03017838602560146
